$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 723.1111
$ws.Range("I33").Value = 881.5714
$ws.Range("K33").Value = 881.5714
$ws.Range("M33").Value = -652.5714
$ws.Range("H43").Value = 6188.0835
$ws.Range("I43").Value = 5361.8335
$ws.Range("K43").Value = 5361.8335
$ws.Range("M43").Value = -5292.8335
$ws.Range("H80").Value = 371.6875
$ws.Range("I80").Value = 367.1
$ws.Range("J80").Value = 379.33334
$ws.Range("K80").Value = 1101.3
$ws.Range("L80").Value = 1138.00002
$ws.Range("M80").Value = -103.3000000000002
$ws.Range("N80").Value = -3134.00002
$ws.Range("H83").Value = 371.6875
$ws.Range("I83").Value = 367.1
$ws.Range("J83").Value = 379.33334
$ws.Range("K83").Value = 3303.9
$ws.Range("L83").Value = 3414.00006
$ws.Range("M83").Value = 1688.1
$ws.Range("N83").Value = -13398.00006
$ws.Range("H107").Value = 536.2083
$ws.Range("I107").Value = 568.86365
$ws.Range("J107").Value = 177
$ws.Range("K107").Value = 568.86365
$ws.Range("L107").Value = 177
$ws.Range("M107").Value = 1351.13635
$ws.Range("N107").Value = -4017
$ws.Range("H111").Value = 2313.2
$ws.Range("I111").Value = 2298
$ws.Range("K111").Value = 6894
$ws.Range("M111").Value = -3827
$ws.Range("H129").Value = 1066.4445
$ws.Range("I129").Value = 1074.75
$ws.Range("K129").Value = 3224.25
$ws.Range("M129").Value = 1775.75
$ws.Range("H132").Value = 1324.3
$ws.Range("I132").Value = 1133.1786
$ws.Range("K132").Value = 3399.5358
$ws.Range("M132").Value = -869.5357999999997
$ws.Range("H137").Value = 1569.6
$ws.Range("I137").Value = 1274.5
$ws.Range("K137").Value = 3823.5
$ws.Range("M137").Value = -1273.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1317.8955
$ws.Range("I32").Value = 1357.0469
$ws.Range("K32").Value = 1357.0469
$ws.Range("M32").Value = -1070.0469
$ws.Range("H34").Value = 99999
$ws.Range("J34").Value = 99999
$ws.Range("L34").Value = 99999
$ws.Range("N34").Value = -100541
$ws.Range("H110").Value = 8440.058999999999
$ws.Range("J110").Value = 9998.5
$ws.Range("L110").Value = 9998.5
$ws.Range("N110").Value = -14088.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 89811.5
$ws.Range("J95").Value = 89811.5
$ws.Range("L95").Value = 89811.5
$ws.Range("N95").Value = -95303.5
$ws.Range("H105").Value = 1701.8
$ws.Range("I105").Value = 1701.8
$ws.Range("K105").Value = 1701.8
$ws.Range("M105").Value = 45.20000000000005
$ws.Range("H107").Value = 2153.3845
$ws.Range("I107").Value = 1844.909
$ws.Range("K107").Value = 1844.909
$ws.Range("M107").Value = 75.09099999999989
$ws.Range("H117").Value = 38924.5
$ws.Range("J117").Value = 38924.5
$ws.Range("L117").Value = 38924.5
$ws.Range("N117").Value = -48102.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3012.976
$ws.Range("I31").Value = 2103.718
$ws.Range("K31").Value = 2103.718
$ws.Range("M31").Value = -1808.718
$ws.Range("H34").Value = 3012.976
$ws.Range("I34").Value = 2103.718
$ws.Range("K34").Value = 2103.718
$ws.Range("M34").Value = -1901.718
$ws.Range("H94").Value = 1055.5714
$ws.Range("J94").Value = 1178
$ws.Range("L94").Value = 1178
$ws.Range("N94").Value = -2080
$ws.Range("H122").Value = 3146.3462
$ws.Range("I122").Value = 2642.5
$ws.Range("J122").Value = 5917.5
$ws.Range("K122").Value = 7927.5
$ws.Range("L122").Value = 17752.5
$ws.Range("M122").Value = -5477.5
$ws.Range("N122").Value = -22652.5
$ws.Range("H141").Value = 34475.777
$ws.Range("I141").Value = 33822.5
$ws.Range("J141").Value = 34998.4
$ws.Range("K141").Value = 33822.5
$ws.Range("L141").Value = 34998.4
$ws.Range("M141").Value = -28642.5
$ws.Range("N141").Value = -45358.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 363.8
$ws.Range("I34").Value = 273
$ws.Range("K34").Value = 819
$ws.Range("M34").Value = -735
$ws.Range("H38").Value = 619
$ws.Range("I38").Value = 37.25
$ws.Range("J38").Value = 798
$ws.Range("K38").Value = 111.75
$ws.Range("L38").Value = 2394
$ws.Range("M38").Value = 235.25
$ws.Range("N38").Value = -3088
$ws.Range("H68").Value = 1011.625
$ws.Range("J68").Value = 959.6
$ws.Range("L68").Value = 2878.8
$ws.Range("N68").Value = -4500.8
$ws.Range("H71").Value = 1011.625
$ws.Range("J71").Value = 959.6
$ws.Range("L71").Value = 8636.4
$ws.Range("N71").Value = -16748.4
$ws.Range("H92").Value = 175
$ws.Range("J92").Value = 175
$ws.Range("L92").Value = 525
$ws.Range("N92").Value = -3021
$ws.Range("H117").Value = 1251761.6
$ws.Range("I117").Value = 900
$ws.Range("J117").Value = 1430456.1
$ws.Range("K117").Value = 2700
$ws.Range("L117").Value = 4291368.300000001
$ws.Range("M117").Value = 742
$ws.Range("N117").Value = -4298252.300000001
$ws.Range("H129").Value = 18126926
$ws.Range("I129").Value = 41793090
$ws.Range("J129").Value = 915171.4399999999
$ws.Range("K129").Value = 125379270
$ws.Range("L129").Value = 2745514.32
$ws.Range("M129").Value = -125374270
$ws.Range("N129").Value = -2755514.32
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3225.037
$ws.Range("I80").Value = 3862.182
$ws.Range("K80").Value = 3862.182
$ws.Range("M80").Value = -2864.182
$ws.Range("H83").Value = 3225.037
$ws.Range("I83").Value = 3862.182
$ws.Range("K83").Value = 19310.91
$ws.Range("M83").Value = -14318.91
$ws.Range("H92").Value = 11375
$ws.Range("J92").Value = 11375
$ws.Range("L92").Value = 11375
$ws.Range("N92").Value = -15119
$ws.Range("H122").Value = 8299.5
$ws.Range("I122").Value = 8169.6
$ws.Range("K122").Value = 24508.8
$ws.Range("M122").Value = -22058.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1049.5
$ws.Range("I40").Value = 899.3333
$ws.Range("K40").Value = 899.3333
$ws.Range("M40").Value = -763.3333
$ws.Range("H132").Value = 8531.833000000001
$ws.Range("I132").Value = 8657.941999999999
$ws.Range("K132").Value = 25973.826
$ws.Range("M132").Value = -23443.826
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 566.8333
$ws.Range("I81").Value = 539.8
$ws.Range("J81").Value = 702
$ws.Range("K81").Value = 1079.6
$ws.Range("L81").Value = 1404
$ws.Range("M81").Value = -18.59999999999991
$ws.Range("N81").Value = -3526
$ws.Range("H84").Value = 566.8333
$ws.Range("I84").Value = 539.8
$ws.Range("J84").Value = 702
$ws.Range("K84").Value = 5398
$ws.Range("L84").Value = 7020
$ws.Range("M84").Value = -94
$ws.Range("N84").Value = -17628
$ws.Range("H107").Value = 462.5
$ws.Range("I107").Value = 462.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1387.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 532.5
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 2296.0833
$ws.Range("I122").Value = 1428.9333
$ws.Range("J122").Value = 3741.3333
$ws.Range("K122").Value = 4286.7999
$ws.Range("L122").Value = 11223.9999
$ws.Range("M122").Value = -1836.7999
$ws.Range("N122").Value = -16123.9999
$ws.Range("H126").Value = 4056.2964
$ws.Range("I126").Value = 3691.3809
$ws.Range("K126").Value = 11074.1427
$ws.Range("M126").Value = -8604.1427
